$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: Update title (column D)
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 36: Fix typo in title (column D)
$ws.Range("D36").Value = "Introduction to unsupervised domain adapation"

# Row 51: Update title (column D) and link (column E)
$ws.Range("D51").Value = "[python] 010-1234-5678과 같은 전화번호 문자열에서 하이픈(-) 빼기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-010-1234-5678%EA%B3%BC-%EA%B0%99%EC%9D%80-%EC%A0%84%ED%99%94%EB%B2%88%ED%98%B8-%EB%AC%B8%EC%9E%90%EC%97%B4%EC%97%90%EC%84%9C-%ED%95%98%EC%9D%B4%ED%94%88-%EB%B9%BC%EA%B8%B0"
